$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Khalid Al Ameri"
$wsSummary.Range("B4").Value = 5923.84
$wsSummary.Range("B6").Value = 472385
$wsSummary.Range("B7").Value = 321072
$wsSummary.Range("B8").Value = 151313
$wsSummary.Range("B9").Value = 1.47

# ---------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Remove the old row 4 (Liquid Assets / Savings Account / 3216);
# TOTAL ASSETS row shifts up from row 5 to row 4.
$wsAssets.Rows.Item(4).Delete()

# Row 2: Vehicles / Mid-range Car / 129601 -> Vehicles / Luxury Car / 463536
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 463536

# Row 3: Vehicles / Mid-range Car / 80111 -> Liquid Assets / Savings Account / 8849
$wsAssets.Range("A3").Value = "Liquid Assets"
$wsAssets.Range("B3").Value = "Savings Account"
$wsAssets.Range("C3").Value = 8849

# Row 4 (formerly row 5): TOTAL ASSETS
$wsAssets.Range("C4").Value = 472385

# ---------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------
$wsLiab = $wb.Worksheets.Item("Liabilities")

# Remove old rows 4 (Personal Loans) and 5 (Credit Cards);
# TOTAL LIABILITIES row shifts up from row 6 to row 4.
$wsLiab.Range("A4:E5").EntireRow.Delete()

# Row 2: Auto Loans / Vehicle Loan 1
$wsLiab.Range("C2").Value = 278122
$wsLiab.Range("D2").Value = 5794
$wsLiab.Range("E2").Value = 4

# Row 3: Auto Loans / Vehicle Loan 2 -> Credit Cards / Credit Card Balance
$wsLiab.Range("A3").Value = "Credit Cards"
$wsLiab.Range("B3").Value = "Credit Card Balance"
$wsLiab.Range("C3").Value = 42950
$wsLiab.Range("D3").Value = 2148
$wsLiab.Range("E3").Value = 1

# Row 4 (formerly row 6): TOTAL LIABILITIES
$wsLiab.Range("C4").Value = 321072

Write-Host "edits applied"
